# ==========================================================================
# #5: cash & deposit done
# Rebuilds the "存款" (deposit) sheet (sheet index 3) to match the finished
# schema used by every other sheet in the workbook: a real header row plus
# the trailing property_category/category/date/legislator_name/legislator_id/
# source_file/index metadata columns. The old sheet had its header row wrongly
# duplicating the first data row, an unused/misplaced amount in column F, and
# the running total stuck in column G with nothing after it - this fixes all
# of that.
# ==========================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

$xlPasteFormats = -4122

# --- Header row (row 1): replace the accidental copy of row 2s data with
#     the real field names, and extend the header across the new columns ---
$ws.Range("F1").Copy()
$ws.Range("H1:M1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# --- Data rows (2-18): move the deposit total from column G into the new
#     "total" column F, then fill in the constant per-record metadata that
#     the other sheets already carry (property_category .. index) ---

# Column I holds literal "2013-11-08" text dates - force text format first
# so Excel does not reinterpret the string as a date serial number.
$ws.Range("I2:I18").NumberFormat = "@"

# record #67 (row 2)
$ws.Range("F2").Copy()
$ws.Range("H2:M2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F2").Value = 124404
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
$ws.Range("I2").Value = "2013-11-08"
$ws.Range("J2").Value = "葉宜津"
$ws.Range("K2").Value = 855
$ws.Range("L2").Value = "tmpabd41"
$ws.Range("M2").Value = 67

# record #68 (row 3)
$ws.Range("F3").Copy()
$ws.Range("H3:M3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F3").Value = 4874
$ws.Range("G3").Value = "deposit"
$ws.Range("H3").Value = "normal"
$ws.Range("I3").Value = "2013-11-08"
$ws.Range("J3").Value = "葉宜津"
$ws.Range("K3").Value = 855
$ws.Range("L3").Value = "tmpabd41"
$ws.Range("M3").Value = 68

# record #69 (row 4)
$ws.Range("F4").Copy()
$ws.Range("H4:M4").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F4").Value = 15896
$ws.Range("G4").Value = "deposit"
$ws.Range("H4").Value = "normal"
$ws.Range("I4").Value = "2013-11-08"
$ws.Range("J4").Value = "葉宜津"
$ws.Range("K4").Value = 855
$ws.Range("L4").Value = "tmpabd41"
$ws.Range("M4").Value = 69

# record #70 (row 5)
$ws.Range("F5").Copy()
$ws.Range("H5:M5").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F5").Value = 12
$ws.Range("G5").Value = "deposit"
$ws.Range("H5").Value = "normal"
$ws.Range("I5").Value = "2013-11-08"
$ws.Range("J5").Value = "葉宜津"
$ws.Range("K5").Value = 855
$ws.Range("L5").Value = "tmpabd41"
$ws.Range("M5").Value = 70

# record #71 (row 6)
$ws.Range("F6").Copy()
$ws.Range("H6:M6").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F6").Value = 1420102
$ws.Range("G6").Value = "deposit"
$ws.Range("H6").Value = "normal"
$ws.Range("I6").Value = "2013-11-08"
$ws.Range("J6").Value = "葉宜津"
$ws.Range("K6").Value = 855
$ws.Range("L6").Value = "tmpabd41"
$ws.Range("M6").Value = 71

# record #72 (row 7)
$ws.Range("F7").Copy()
$ws.Range("H7:M7").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F7").Value = 391966
$ws.Range("G7").Value = "deposit"
$ws.Range("H7").Value = "normal"
$ws.Range("I7").Value = "2013-11-08"
$ws.Range("J7").Value = "葉宜津"
$ws.Range("K7").Value = 855
$ws.Range("L7").Value = "tmpabd41"
$ws.Range("M7").Value = 72

# record #73 (row 8)
$ws.Range("F8").Copy()
$ws.Range("H8:M8").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F8").Value = 50821
$ws.Range("G8").Value = "deposit"
$ws.Range("H8").Value = "normal"
$ws.Range("I8").Value = "2013-11-08"
$ws.Range("J8").Value = "葉宜津"
$ws.Range("K8").Value = 855
$ws.Range("L8").Value = "tmpabd41"
$ws.Range("M8").Value = 73

# record #74 (row 9)
$ws.Range("F9").Copy()
$ws.Range("H9:M9").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F9").Value = 428610
$ws.Range("G9").Value = "deposit"
$ws.Range("H9").Value = "normal"
$ws.Range("I9").Value = "2013-11-08"
$ws.Range("J9").Value = "葉宜津"
$ws.Range("K9").Value = 855
$ws.Range("L9").Value = "tmpabd41"
$ws.Range("M9").Value = 74

# record #75 (row 10)
$ws.Range("F10").Copy()
$ws.Range("H10:M10").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F10").Value = 109101
$ws.Range("G10").Value = "deposit"
$ws.Range("H10").Value = "normal"
$ws.Range("I10").Value = "2013-11-08"
$ws.Range("J10").Value = "葉宜津"
$ws.Range("K10").Value = 855
$ws.Range("L10").Value = "tmpabd41"
$ws.Range("M10").Value = 75

# record #76 (row 11)
$ws.Range("F11").Copy()
$ws.Range("H11:M11").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F11").Value = 104778
$ws.Range("G11").Value = "deposit"
$ws.Range("H11").Value = "normal"
$ws.Range("I11").Value = "2013-11-08"
$ws.Range("J11").Value = "葉宜津"
$ws.Range("K11").Value = 855
$ws.Range("L11").Value = "tmpabd41"
$ws.Range("M11").Value = 76

# record #77 (row 12)
$ws.Range("F12").Copy()
$ws.Range("H12:M12").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F12").Value = 228003
$ws.Range("G12").Value = "deposit"
$ws.Range("H12").Value = "normal"
$ws.Range("I12").Value = "2013-11-08"
$ws.Range("J12").Value = "葉宜津"
$ws.Range("K12").Value = 855
$ws.Range("L12").Value = "tmpabd41"
$ws.Range("M12").Value = 77

# record #78 (row 13)
$ws.Range("F13").Copy()
$ws.Range("H13:M13").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F13").Value = 479391.85
$ws.Range("G13").Value = "deposit"
$ws.Range("H13").Value = "normal"
$ws.Range("I13").Value = "2013-11-08"
$ws.Range("J13").Value = "葉宜津"
$ws.Range("K13").Value = 855
$ws.Range("L13").Value = "tmpabd41"
$ws.Range("M13").Value = 78

# record #79 (row 14)
$ws.Range("F14").Copy()
$ws.Range("H14:M14").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F14").Value = 83084
$ws.Range("G14").Value = "deposit"
$ws.Range("H14").Value = "normal"
$ws.Range("I14").Value = "2013-11-08"
$ws.Range("J14").Value = "葉宜津"
$ws.Range("K14").Value = 855
$ws.Range("L14").Value = "tmpabd41"
$ws.Range("M14").Value = 79

# record #80 (row 15)
$ws.Range("F15").Copy()
$ws.Range("H15:M15").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F15").Value = 526996
$ws.Range("G15").Value = "deposit"
$ws.Range("H15").Value = "normal"
$ws.Range("I15").Value = "2013-11-08"
$ws.Range("J15").Value = "葉宜津"
$ws.Range("K15").Value = 855
$ws.Range("L15").Value = "tmpabd41"
$ws.Range("M15").Value = 80

# record #81 (row 16)
$ws.Range("F16").Copy()
$ws.Range("H16:M16").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F16").Value = 522673
$ws.Range("G16").Value = "deposit"
$ws.Range("H16").Value = "normal"
$ws.Range("I16").Value = "2013-11-08"
$ws.Range("J16").Value = "葉宜津"
$ws.Range("K16").Value = 855
$ws.Range("L16").Value = "tmpabd41"
$ws.Range("M16").Value = 81

# record #82 (row 17)
$ws.Range("F17").Copy()
$ws.Range("H17:M17").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F17").Value = 1936
$ws.Range("G17").Value = "deposit"
$ws.Range("H17").Value = "normal"
$ws.Range("I17").Value = "2013-11-08"
$ws.Range("J17").Value = "葉宜津"
$ws.Range("K17").Value = 855
$ws.Range("L17").Value = "tmpabd41"
$ws.Range("M17").Value = 82

# record #83 (row 18)
$ws.Range("F18").Copy()
$ws.Range("H18:M18").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F18").Value = 131506
$ws.Range("G18").Value = "deposit"
$ws.Range("H18").Value = "normal"
$ws.Range("I18").Value = "2013-11-08"
$ws.Range("J18").Value = "葉宜津"
$ws.Range("K18").Value = 855
$ws.Range("L18").Value = "tmpabd41"
$ws.Range("M18").Value = 83

